# Auto-generated edit script
# Applies numeric updates to the Leve profit-tracking tables (H-N columns)
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching
# a scheduled market-price refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2391.3928
$ws.Range("I19").Value = 811.8
$ws.Range("J19").Value = 3268.9443
$ws.Range("K19").Value = 811.8
$ws.Range("L19").Value = 3268.9443
$ws.Range("M19").Value = -636.8
$ws.Range("N19").Value = -3618.9443

$ws.Range("H70").Value = 864390
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 864390
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 2593170
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -2593710

$ws.Range("H73").Value = 864390
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 864390
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 2593170
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -2595042

$ws.Range("H74").Value = 7637.857
$ws.Range("I74").Value = 7485.5
$ws.Range("K74").Value = 7485.5
$ws.Range("M74").Value = -6549.5

$ws.Range("H77").Value = 7637.857
$ws.Range("I77").Value = 7485.5
$ws.Range("K77").Value = 37427.5
$ws.Range("M77").Value = -32747.5

$ws.Range("H132").Value = 1554.25
$ws.Range("I132").Value = 1604.1818
$ws.Range("K132").Value = 4812.5454
$ws.Range("M132").Value = -2282.5454

$ws.Range("H134").Value = 120000
$ws.Range("J134").Value = 120000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -130140

$ws.Range("H137").Value = 3221.65
$ws.Range("J137").Value = 4420.7646
$ws.Range("L137").Value = 13262.2938
$ws.Range("N137").Value = -18362.2938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28586
$ws.Range("I32").Value = 28937.588
$ws.Range("K32").Value = 28937.588
$ws.Range("M32").Value = -28650.588

$ws.Range("H74").Value = 6200.7744
$ws.Range("I74").Value = 5005.5835
$ws.Range("K74").Value = 5005.5835
$ws.Range("M74").Value = -4131.5835

$ws.Range("H77").Value = 6200.7744
$ws.Range("I77").Value = 5005.5835
$ws.Range("K77").Value = 25027.9175
$ws.Range("M77").Value = -20659.9175

$ws.Range("H88").Value = 1440.9166
$ws.Range("I88").Value = 1213.2858
$ws.Range("J88").Value = 1759.6
$ws.Range("K88").Value = 1213.2858
$ws.Range("L88").Value = 1759.6
$ws.Range("M88").Value = -807.2858000000001
$ws.Range("N88").Value = -2571.6

$ws.Range("H91").Value = 1440.9166
$ws.Range("I91").Value = 1213.2858
$ws.Range("J91").Value = 1759.6
$ws.Range("K91").Value = 1213.2858
$ws.Range("L91").Value = 1759.6
$ws.Range("M91").Value = 190.7141999999999
$ws.Range("N91").Value = -4567.6

$ws.Range("H92").Value = 55000
$ws.Range("J92").Value = 55000
$ws.Range("L92").Value = 55000
$ws.Range("N92").Value = -59992

$ws.Range("H94").Value = 20499.5
$ws.Range("J94").Value = 20499.5
$ws.Range("L94").Value = 20499.5
$ws.Range("N94").Value = -22301.5

$ws.Range("H110").Value = 11365479
$ws.Range("I110").Value = 22728096
$ws.Range("J110").Value = 2862.7273
$ws.Range("K110").Value = 22728096
$ws.Range("L110").Value = 2862.7273
$ws.Range("M110").Value = -22726051
$ws.Range("N110").Value = -6952.7273

$ws.Range("H128").Value = 210624.25
$ws.Range("J128").Value = 210624.25
$ws.Range("L128").Value = 210624.25
$ws.Range("N128").Value = -220584.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2822.1538
$ws.Range("I20").Value = 2102.5264
$ws.Range("K20").Value = 2102.5264
$ws.Range("M20").Value = -1855.5264

$ws.Range("H94").Value = 2229.6875
$ws.Range("I94").Value = 2518.6667
$ws.Range("J94").Value = 1858.1428
$ws.Range("K94").Value = 2518.6667
$ws.Range("L94").Value = 1858.1428
$ws.Range("M94").Value = -2067.6667
$ws.Range("N94").Value = -2760.1428

$ws.Range("H120").Value = 66336.664
$ws.Range("J120").Value = 66336.664
$ws.Range("L120").Value = 66336.664
$ws.Range("N120").Value = -76012.664

$ws.Range("H124").Value = 46780
$ws.Range("J124").Value = 46780
$ws.Range("L124").Value = 46780
$ws.Range("N124").Value = -56600

$ws.Range("H125").Value = 38389.75
$ws.Range("J125").Value = 38389.75
$ws.Range("L125").Value = 38389.75
$ws.Range("N125").Value = -48229.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3736.3438
$ws.Range("I31").Value = 2093.05
$ws.Range("J31").Value = 6475.1665
$ws.Range("K31").Value = 2093.05
$ws.Range("L31").Value = 6475.1665
$ws.Range("M31").Value = -1798.05
$ws.Range("N31").Value = -7065.1665

$ws.Range("H34").Value = 3736.3438
$ws.Range("I34").Value = 2093.05
$ws.Range("J34").Value = 6475.1665
$ws.Range("K34").Value = 2093.05
$ws.Range("L34").Value = 6475.1665
$ws.Range("M34").Value = -1891.05
$ws.Range("N34").Value = -6879.1665

$ws.Range("H60").Value = 30000
$ws.Range("J60").Value = 30000
$ws.Range("L60").Value = 30000
$ws.Range("N60").Value = -31022

$ws.Range("H105").Value = 1480.5264
$ws.Range("I105").Value = 1480.5264
$ws.Range("K105").Value = 1480.5264
$ws.Range("M105").Value = 266.4736

$ws.Range("H134").Value = 3798.8445
$ws.Range("I134").Value = 3309.95
$ws.Range("K134").Value = 9929.849999999999
$ws.Range("M134").Value = -7394.849999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 3666.6667
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 3666.6667
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H88").Value = 5414.2856
$ws.Range("I88").Value = 4650
$ws.Range("J88").Value = 10000
$ws.Range("K88").Value = 13950
$ws.Range("L88").Value = 30000
$ws.Range("M88").Value = -13522
$ws.Range("N88").Value = -30856

$ws.Range("H91").Value = 5414.2856
$ws.Range("I91").Value = 4650
$ws.Range("J91").Value = 10000
$ws.Range("K91").Value = 13950
$ws.Range("L91").Value = 30000
$ws.Range("M91").Value = -12468
$ws.Range("N91").Value = -32964

$ws.Range("H140").Value = 2103.2
$ws.Range("I140").Value = 1276.2858
$ws.Range("J140").Value = 4032.6667
$ws.Range("K140").Value = 3828.8574
$ws.Range("L140").Value = 12098.0001
$ws.Range("M140").Value = 1351.1426
$ws.Range("N140").Value = -22458.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4090
$ws.Range("I43").Value = 5385
$ws.Range("J43").Value = 1500
$ws.Range("K43").Value = 5385
$ws.Range("L43").Value = 1500
$ws.Range("M43").Value = -5234
$ws.Range("N43").Value = -1802

$ws.Range("H70").Value = 6339.5
$ws.Range("I70").Value = 7947.25
$ws.Range("J70").Value = 4731.75
$ws.Range("K70").Value = 7947.25
$ws.Range("L70").Value = 4731.75
$ws.Range("M70").Value = -7677.25
$ws.Range("N70").Value = -5271.75

$ws.Range("H73").Value = 6339.5
$ws.Range("I73").Value = 7947.25
$ws.Range("J73").Value = 4731.75
$ws.Range("K73").Value = 7947.25
$ws.Range("L73").Value = 4731.75
$ws.Range("M73").Value = -7011.25
$ws.Range("N73").Value = -6603.75

$ws.Range("H132").Value = 5975.7856
$ws.Range("I132").Value = 4160.2905
$ws.Range("K132").Value = 12480.8715
$ws.Range("M132").Value = -9950.871500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4094.3333
$ws.Range("I68").Value = 3066.5
$ws.Range("K68").Value = 3066.5
$ws.Range("M68").Value = -2317.5

$ws.Range("H71").Value = 4094.3333
$ws.Range("I71").Value = 3066.5
$ws.Range("K71").Value = 15332.5
$ws.Range("M71").Value = -11588.5

$ws.Range("H82").Value = 1877.1111
$ws.Range("I82").Value = 1749
$ws.Range("K82").Value = 1749
$ws.Range("M82").Value = -1388

$ws.Range("H85").Value = 1877.1111
$ws.Range("I85").Value = 1749
$ws.Range("K85").Value = 1749
$ws.Range("M85").Value = -501

$ws.Range("H93").Value = 2045.8334
$ws.Range("I93").Value = 2333.3333
$ws.Range("J93").Value = 1758.3334
$ws.Range("K93").Value = 2333.3333
$ws.Range("L93").Value = 1758.3334
$ws.Range("M93").Value = -1085.3333
$ws.Range("N93").Value = -4254.3334

$ws.Range("H100").Value = 15628587
$ws.Range("I100").Value = 41669784
$ws.Range("K100").Value = 41669784
$ws.Range("M100").Value = -41669243

$ws.Range("H124").Value = 77351.2
$ws.Range("J124").Value = 77351.2
$ws.Range("L124").Value = 77351.2
$ws.Range("N124").Value = -87171.2

$ws.Range("H136").Value = 1964255.1
$ws.Range("I136").Value = 2566410.8
$ws.Range("K136").Value = 7699232.399999999
$ws.Range("M136").Value = -7696682.399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12812.919
$ws.Range("J81").Value = 16882.521
$ws.Range("L81").Value = 33765.042
$ws.Range("N81").Value = -35887.042

$ws.Range("H84").Value = 12812.919
$ws.Range("J84").Value = 16882.521
$ws.Range("L84").Value = 168825.21
$ws.Range("N84").Value = -179433.21

$ws.Range("H96").Value = 3732.75
$ws.Range("I96").Value = 3595.4
$ws.Range("J96").Value = 3961.6667
$ws.Range("K96").Value = 3595.4
$ws.Range("L96").Value = 3961.6667
$ws.Range("M96").Value = -2222.4
$ws.Range("N96").Value = -6707.6667

$ws.Range("H120").Value = 50460
$ws.Range("J120").Value = 50460
$ws.Range("L120").Value = 50460
$ws.Range("N120").Value = -60136

$ws.Range("H136").Value = 4204240
$ws.Range("I136").Value = 5495645
$ws.Range("J136").Value = 7173.125
$ws.Range("K136").Value = 16486935
$ws.Range("L136").Value = 21519.375
$ws.Range("M136").Value = -16484385
$ws.Range("N136").Value = -26619.375
